$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear J2 and K2
$ws.Range("J2").Value = $null
$ws.Range("K2").Value = $null

# Row 3: J3 -> 5, K3 -> blank, N3 -> 5
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = $null
$ws.Range("N3").Value = 5

# Row 4: J4 -> 5, K4 -> 3, N4 -> 5
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 3
$ws.Range("N4").Value = 5

# Row 5: J5 -> 5, N5 -> 5
$ws.Range("J5").Value = 5
$ws.Range("N5").Value = 5

# Row 6: J6 -> 5, K6 -> blank, N6 -> 3
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = $null
$ws.Range("N6").Value = 3
